$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 45607
$ws.Range("A5").NumberFormat = "d-mmm"
$ws.Range("B5").Value = "Biwak"
$ws.Range("C5").Value = "Szkoła"

# Row 6
$ws.Range("A6").Value = 45547
$ws.Range("A6").NumberFormat = "d-mmm"
$ws.Range("C6").Value = "Oleszna"
$ws.Range("D6").Value = "Przemek"

# Row 7
$ws.Range("B7").Value = "Nic "
$ws.Range("C7").Value = "Nigdzie"
$ws.Range("D7").Value = "Nikt"

# Update selection to match the final active cell
$ws.Range("D7").Select() | Out-Null
